$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-12-04 12:52:59"

for ($r = 2; $r -le 18; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
